$d = $word.ActiveDocument

$replacements = @(
    @("987÷8=", "559÷2="),
    @("437÷4=", "564÷7="),
    @("939÷2=", "726÷6="),
    @("593÷2=", "407÷9="),
    @("423÷9=", "896÷8="),
    @("526÷8=", "811÷9="),
    @("537÷3=", "572÷7="),
    @("634÷3=", "226÷3="),
    @("542÷8=", "635÷6="),
    @("278÷4=", "768÷6="),
    @("536÷8=", "297÷2="),
    @("524÷2=", "460÷8="),
    @("933÷8=", "740÷5="),
    @("744÷9=", "595÷9="),
    @("250÷9=", "561÷7="),
    @("303÷7=", "504÷6="),
    @("376÷3=", "189÷9="),
    @("920÷5=", "875÷2="),
    @("861÷4=", "604÷8="),
    @("112÷2=", "532÷4="),
    @("878÷6=", "277÷9="),
    @("133÷8=", "294÷9="),
    @("339÷3=", "728÷6="),
    @("106÷8=", "427÷3="),
    @("939÷9=", "219÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
